$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 97, pushing the existing row 97 (and below) down to 98.
$ws.Rows.Item(97).Insert()

# New row 97: same market / region as the row that used to be here (now row 98),
# but with the new date + price data.
$ws.Cells.Item(97, 1).Value = 1
$ws.Cells.Item(97, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(97, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(97, 4).Value = 45267
$ws.Cells.Item(97, 5).Value = 15
$ws.Cells.Item(97, 6).Value = 100112028
$ws.Cells.Item(97, 7).Value = "Sandia"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 850
$ws.Cells.Item(97, 11).Value = 480
$ws.Cells.Item(97, 12).Value = 500
$ws.Cells.Item(97, 13).Value = 488
$ws.Cells.Item(97, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(97, 15).Value = "Perú"
$ws.Cells.Item(97, 16).Value = 488
$ws.Cells.Item(97, 17).Value = 1
$ws.Cells.Item(97, 18).Value = "Hortaliza"
